# Expense details: insert a new "Headphones" expense row above the existing
# "Rent" row (old row 2 data shifts down to row 3).
#
# Before:
#   Row1: Category | Amount | Data
#   Row2: (blank)  | 300    | 45688.87447916667   (date-formatted)
#
# After:
#   Row1: Category   | Amount | Data
#   Row2: Headphones | 75     | 45856.8328125      (date-formatted)
#   Row3: Rent       | 300    | 45688.87447916667  (date-formatted)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing B2:C2 (amount + formatted date) down into row 3 first,
# carrying their number formatting (style) along with them, so the old
# "Rent" data keeps its date format in its new home.
$ws.Range("B2:C2").Copy($ws.Range("B3:C3"))
$ws.Range("A3").Value = "Rent"

# Now populate the freed-up row 2 with the new "Headphones" expense.
$ws.Range("A2").Value = "Headphones"
$ws.Range("B2").Value = 75
$ws.Range("C2").Value = 45856.8328125
